$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("app.controller(""main_ctrl"", function(`$scope,`$http) {")
Write-Output $find.Found
